# gym_log_Q1_2024 - bio_data.xlsx
# "Code update, data source update"
# Appends the new bio-tracking rows (2024-06-14 .. 2024-06-26) that were
# logged after the last save, mirroring the existing Date/Waist/Weight/
# kcal/kcal total/Creatine columns on List1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List1")

# Date, Waist, Weight, kcal, kcal total, Creatine
$newData = @(
    @("2024-06-14", 97,    80.9,  2680, 2680, $null),
    @("2024-06-15", 97,    80.1,  2554, 2554, $null),
    @("2024-06-16", 97,    81,    2942, 2942, $null),
    @("2024-06-17", 97.5,  81.1,  2910, 2910, $null),
    @("2024-06-18", 98,    81.4,  2389, 2277, $null),
    @("2024-06-19", 98.5,  81.9,  2438, 2303, $null),
    @("2024-06-20", 98.5,  81.5,  1790, 1621, $null),
    @("2024-06-21", 97.5,  80.2,  2190, 2190, $null),
    @("2024-06-22", 97.5,  80.6,  3494, 3494, $null),
    @("2024-06-23", 98,    82,    2755, 2755, $null),
    @("2024-06-24", 98,    81,    2111, 1943, 1),
    @("2024-06-25", 97,    79.8,  1753, 1555, 1),
    @("2024-06-26", 97,    79.7,  2031, 2031, 1)
)

$startRow = 162
for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $entry = $newData[$i]

    $ws.Range("A$r").Value = $entry[0]
    $ws.Range("B$r").Value = $entry[1]
    $ws.Range("C$r").Value = $entry[2]
    $ws.Range("D$r").Value = $entry[3]
    $ws.Range("E$r").Value = $entry[4]
    if ($entry[5] -ne $null) {
        $ws.Range("F$r").Value = $entry[5]
    }
}

# Match the author's final scroll position / active cell from the commit.
$ws.Range("C170").Select()
